$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, pushing existing rows 6..86 down to 7..87.
$ws.Rows.Item(6).Insert()

# Fill in the new row 6 with its data.
$ws.Range("A6").Value = 10
$ws.Range("B6").Value = "Vega Modelo de Temuco"
$ws.Range("C6").Value = "La Araucanía"
$ws.Range("D6").Value = (Get-Date -Year 2022 -Month 11 -Day 17 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E6").Value = 9
$ws.Range("F6").Value = 100112026
$ws.Range("G6").Value = "Haba"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 180
$ws.Range("K6").Value = 10000
$ws.Range("L6").Value = 11000
$ws.Range("M6").Value = 10444
$ws.Range("N6").Value = "$/saco 25 kilos"
$ws.Range("O6").Value = "Región Metropolitana"
$ws.Range("P6").Value = 418
$ws.Range("Q6").Value = 25
$ws.Range("R6").Value = "Hortaliza"
